$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: "20min" -> "20 min" ---
$ws.Range("B16").Value = "20 min"

# --- Row 17: new entry ---
$ws.Range("A17").Value = "fix detail page"
$ws.Range("B17").Value = "5 min"
$ws.Range("C17").Value = 45926
$ws.Range("D17").Value = "fix: increase media query breakpoint for mobile responsiveness on house details page"
$ws.Range("D17").WrapText = $true

# --- Row 18: new entry ---
$ws.Range("A18").Value = "fix console.log"
$ws.Range("D18").Value = "fix: remove all console.logs"
$ws.Range("B18").Value = "3 min"
$ws.Range("C18").Value = 45926

# --- Row 19: new entry (tall "feature" row) ---
$ws.Range("A19").Value = "refactor components/views"
$ws.Range("B19").Value = "1:30 h"
$ws.Range("C19").Value = 45926
$ws.Range("D19").Value = "- Introduced FormRow and FormField components for better form structure and reusability. - Implemented ImageUpload component for handling image uploads with preview functionality. - Created SelectField component for dropdown selections. - Replaced existing form elements in CreateHouseView with new components for consistency and maintainability. - Added NoResults component to standardize no results display across views. - Refactored FavoritesView and HousesView to utilize new components for search and sorting functionality."
$ws.Range("A19").VerticalAlignment = -4108
$ws.Range("B19").VerticalAlignment = -4108
$ws.Range("C19").VerticalAlignment = -4108
$ws.Range("D19").VerticalAlignment = -4108
$ws.Range("D19").WrapText = $true
$ws.Rows("19").RowHeight = 213

# --- Row 20: new entry ---
$ws.Range("A20").Value = "fix form errors"
$ws.Range("B20").Value = "15 min"
$ws.Range("C20").Value = 45926
$ws.Range("D20").Value = "enhance form error handeling by adding red border and red place holder"
$ws.Range("C20").VerticalAlignment = -4108
$ws.Range("D20").Borders.Item(8).LineStyle = -4142
$ws.Range("D20").Borders.Item(9).LineStyle = -4142
$ws.Range("D20").Borders.Item(7).Color = 0xAAAAAA
$ws.Range("D20").Borders.Item(10).Color = 0xAAAAAA

# --- Row 30: updated total ---
$ws.Range("B30").Value = "14 hours and 58 minutes"

# --- Misc view state (best-effort match of author's final selection) ---
$ws.Range("I22").Select()
